$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "TestUser_1DKUKxoJrIT"
$ws.Range("B7").Value = "automation_test+1611435001+RuarLPPYUp@gmail.com"
$ws.Range("C7").Value = "TestPassw0rd@123!`$lOOtqQERdv"
